# Update the "Billable" (column E) flag from 1 to 0 for the specified employees
# in the Employees worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")

$rows = @(4, 7, 13, 15, 16, 17, 18, 22, 23, 24, 25, 26, 29, 33, 37, 42, 43, 48, 50, 74, 79)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = 0
}
